# "edits to gifs and extrapolations graph"
#
# - Move the selection on Hoja2 from B22 to B16 (and it is no longer the
#   active/selected tab).
# - Add a new worksheet named "extrapolations" after Hoja2, becoming the
#   new active tab, containing a small weight/time_period/equivalence
#   table, with the selection left on C6.

$wb = $excel.ActiveWorkbook

# --- Update selection on Hoja2 -------------------------------------------
$hoja2 = $wb.Worksheets.Item("Hoja2")
$hoja2.Range("B16").Select()

# --- Add the new "extrapolations" worksheet at the end -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "extrapolations"

# Fill in the table. The order of entry matches the order new shared
# strings were authored in the original workbook.
$ws.Range("A1").Value = "weight"

$ws.Range("C3").Value = "1.5 car engines"
$ws.Range("C4").Value = "1 Toyota"
$ws.Range("C5").Value = "5 Hummers"

$ws.Range("C1").Value = "equivalence"
$ws.Range("B1").Value = "time_period"

$ws.Range("C2").Value = "4 car tires"

$ws.Range("B3").Value = "1 week"
$ws.Range("B4").Value = "1 month"

$ws.Range("B2").Value = "1 day"
$ws.Range("B5").Value = "1 year"

$ws.Range("A2").Value = 80
$ws.Range("A3").Value = 560
$ws.Range("A4").Value = 2240
$ws.Range("A5").Value = 26880

# Leave the selection where the author left it.
$ws.Range("C6").Select()
